$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$NewValue)
    $origStyle = $Cell.Style
    $Cell.NumberFormat = "@"
    $Cell.Value = $NewValue
    $Cell.Style = $origStyle
}

# Row 2
Set-TextValue $ws.Cells.Item(2, 4) "243.37"

# Row 3
Set-TextValue $ws.Cells.Item(3, 4) "23.08"

# Row 4
Set-TextValue $ws.Cells.Item(4, 4) "5.399"

# Row 7
Set-TextValue $ws.Cells.Item(7, 4) "6.498"

# Row 8
Set-TextValue $ws.Cells.Item(8, 4) "0.8093"

# Row 9
Set-TextValue $ws.Cells.Item(9, 4) "0.9268"

# Row 11
Set-TextValue $ws.Cells.Item(11, 4) "0.07406"

# Row 12
Set-TextValue $ws.Cells.Item(12, 4) "0.03289"

# Row 13
Set-TextValue $ws.Cells.Item(13, 4) "0.03063"

# Row 14
Set-TextValue $ws.Cells.Item(14, 4) "0.09353"

# Row 15
Set-TextValue $ws.Cells.Item(15, 4) "3.846"

# Row 16
Set-TextValue $ws.Cells.Item(16, 4) "0.001584"

# Row 17
Set-TextValue $ws.Cells.Item(17, 4) "0.04685"

# Row 18
Set-TextValue $ws.Cells.Item(18, 2) "One"
Set-TextValue $ws.Cells.Item(18, 3) "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws.Cells.Item(18, 4) "0.0005919"
Set-TextValue $ws.Cells.Item(18, 5) "17OneONE"

# Row 19
Set-TextValue $ws.Cells.Item(19, 2) "TigerCash"
Set-TextValue $ws.Cells.Item(19, 3) "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Cells.Item(19, 4) "0.005880"
Set-TextValue $ws.Cells.Item(19, 5) "18TigerCashTCH"

# Row 20
Set-TextValue $ws.Cells.Item(20, 2) "BitKan"
Set-TextValue $ws.Cells.Item(20, 3) "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue $ws.Cells.Item(20, 4) "0.001268"
Set-TextValue $ws.Cells.Item(20, 5) "19BitKanKANBestin24h"

# Row 21
Set-TextValue $ws.Cells.Item(21, 2) "HotbitToken"
Set-TextValue $ws.Cells.Item(21, 3) "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue $ws.Cells.Item(21, 4) "0.004886"
Set-TextValue $ws.Cells.Item(21, 5) "20HotbitTokenHTB"

# Row 22
Set-TextValue $ws.Cells.Item(22, 2) "NitroEx"
Set-TextValue $ws.Cells.Item(22, 3) "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue $ws.Cells.Item(22, 4) "0.00007998"
Set-TextValue $ws.Cells.Item(22, 5) "21NitroExNTX"

# Row 23
Set-TextValue $ws.Cells.Item(23, 2) "LEO"
Set-TextValue $ws.Cells.Item(23, 3) "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Cells.Item(23, 4) "3.570"
Set-TextValue $ws.Cells.Item(23, 5) "22LEOLEO"

# Row 24
Set-TextValue $ws.Cells.Item(24, 2) "BTSEToken"
Set-TextValue $ws.Cells.Item(24, 3) "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws.Cells.Item(24, 4) "2.133"
Set-TextValue $ws.Cells.Item(24, 5) "23BTSETokenBTSE"

# Row 40
Set-TextValue $ws.Cells.Item(40, 4) "0.03964"

# Row 41
Set-TextValue $ws.Cells.Item(41, 2) "BKEXToken"
Set-TextValue $ws.Cells.Item(41, 3) "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Cells.Item(41, 4) "0.1078"
Set-TextValue $ws.Cells.Item(41, 5) "40BKEXTokenBKK"

# Row 42
Set-TextValue $ws.Cells.Item(42, 4) "0.002649"

# Row 43
Set-TextValue $ws.Cells.Item(43, 2) "KickToken"
Set-TextValue $ws.Cells.Item(43, 3) "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws.Cells.Item(43, 4) "0.003068"
Set-TextValue $ws.Cells.Item(43, 5) "42KickTokenKICKWorstin24h"

# Row 44
Set-TextValue $ws.Cells.Item(44, 4) "0.008424"

# Row 45
Set-TextValue $ws.Cells.Item(45, 4) "0.00005091"

# Row 47
Set-TextValue $ws.Cells.Item(47, 5) "46CoinbaseStockTokenCOIN"

# Row 48
Set-TextValue $ws.Cells.Item(48, 4) "0.002265"
